# Replace the VLOOKUP formulas in row 3 (columns G:AM) of the "Answers"
# sheet with their resolved literal numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Answers")

$values = [ordered]@{
    "G3"  = 2
    "H3"  = 1
    "I3"  = 5
    "J3"  = 2
    "K3"  = 2
    "L3"  = 5
    "M3"  = 2
    "N3"  = 2
    "O3"  = 4
    "P3"  = 2
    "Q3"  = 2
    "R3"  = 2
    "S3"  = 1
    "T3"  = 3
    "U3"  = 2
    "V3"  = 2
    "W3"  = 5
    "X3"  = 4
    "Y3"  = 2
    "Z3"  = 4
    "AA3" = 5
    "AB3" = 4
    "AC3" = 2
    "AD3" = 3
    "AE3" = 1
    "AF3" = 3
    "AG3" = 4
    "AH3" = 2
    "AI3" = 4
    "AJ3" = 1
    "AK3" = 3
    "AL3" = 4
    "AM3" = 6
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
